# Updated OutlookLogin and OpenAir
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: add a trailing (mostly blank, hyperlink-styled) row at row 15
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("G15").Style = "Hyperlink"

# Sheet1's selection moves to a full column A:B selection (no single
# active cell highlighted, tab no longer "selected" once Sheet2 is
# activated below).
$ws1.Range("A1:B1048576").Select()

# ---------------------------------------------------------------------
# Sheet2 (3rd tab): populate the previously-empty sheet with the
# OutlookLogin / OpenAir test-data table.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

# New shared-string values must be introduced in this exact order so the
# workbook's shared-string table lines up cell-for-cell with the table
# being typed in (column F first, then C, then the F/D/E pairs, etc.)
$ws2.Range("F1").Value = "lifecycle"
$ws2.Range("C1").Value = "depends_TC"
$ws2.Range("F2").Value = "outlookLogin"
$ws2.Range("D2").Value = "http://mail.maveric-systems.com/"
$ws2.Range("E1").Value = "driverType"
$ws2.Range("E2").Value = "desktop"
$ws2.Range("G2").Value = "danielf@maveric-systems.com"
$ws2.Range("H2").Value = "Mavaug@123"

# Remaining cells reuse already-existing shared strings.
$ws2.Range("A1").Value = "C"
$ws2.Range("B1").Value = "TC_ID"
$ws2.Range("D1").Value = "url"
$ws2.Range("G1").Value = "userName"
$ws2.Range("H1").Value = "password"

$ws2.Range("A2").Value = "Y"
$ws2.Range("B2").Value = "TC001"

$ws2.Range("B3").Value = "TC002"
$ws2.Range("B4").Value = "TC003"
$ws2.Range("B5").Value = "TC004"
$ws2.Range("B6").Value = "TC005"
$ws2.Range("B7").Value = "TC006"

# Size the populated columns to fit their new content (mirrors the
# bestFit column widths Excel stamps after typing the data in).
$ws2.Columns("A").AutoFit()
$ws2.Columns("C").AutoFit()
$ws2.Columns("D").AutoFit()
$ws2.Columns("E").AutoFit()
$ws2.Columns("F").AutoFit()
$ws2.Columns("G").AutoFit()
$ws2.Columns("H").AutoFit()

# Sheet2 becomes the active tab, with G10 selected.
$ws2.Activate()
$ws2.Range("G10").Select()
